# Applies the "added new columns in report, changed log formation" edit to
# the DriverSheet: populates three new data columns (H:J) for rows 2-5 and
# moves the active selection to K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DriverSheet")
$ws.Activate()

# New report columns - TestCaseDescription (H) / Validation (I) / ExpectedValidation (J)
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 111

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 22
$ws.Range("J3").Value = 222

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 33
$ws.Range("J4").Value = 333

$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 44
$ws.Range("J5").Value = 444

# Log formation: move the active cell/selection to K5
$ws.Range("K5").Select()
